$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150, shifting existing rows 150-211 down to 151-212
$ws.Rows("150:150").Insert()

# Populate the new row 150 with data (most columns copied from the old row 150 / neighboring rows,
# a few columns updated with new values per the edit)
$ws.Cells.Item(150, 1).Value = 6
$ws.Cells.Item(150, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(150, 3).Value = "Metropolitana"
$ws.Cells.Item(150, 4).Value = 44704
$ws.Cells.Item(150, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(150, 5).Value = 13
$ws.Cells.Item(150, 6).Value = 100112001
$ws.Cells.Item(150, 7).Value = "Berenjena"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 260
$ws.Cells.Item(150, 11).Value = 5000
$ws.Cells.Item(150, 12).Value = 5000
$ws.Cells.Item(150, 13).Value = 5000
$ws.Cells.Item(150, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(150, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(150, 16).Value = 100
$ws.Cells.Item(150, 17).Value = 50
$ws.Cells.Item(150, 18).Value = "Hortaliza"
